$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date/id columns stay as plain text (matching existing rows) instead
# of being auto-converted to a date serial number / numeric value. Apply a
# text format before assigning, then clear the format again so the new
# cells don't pick up an extra style compared to the rest of the sheet.
$ws.Range("B135:C136").NumberFormat = "@"

# Row 135: 2020-02-27
$ws.Cells.Item(135, 1).Value = 1582761600
$ws.Cells.Item(135, 2).Value = "2020-02-27"
$ws.Cells.Item(135, 3).Value = "0211"
$ws.Cells.Item(135, 4).Value = "TASHIN"
$ws.Cells.Item(135, 5).Value = 0.24
$ws.Cells.Item(135, 6).Value = 0.245
$ws.Cells.Item(135, 7).Value = 0.23
$ws.Cells.Item(135, 8).Value = 0.23
$ws.Cells.Item(135, 9).Value = 334300

# Row 136: 2020-02-28
$ws.Cells.Item(136, 1).Value = 1582848000
$ws.Cells.Item(136, 2).Value = "2020-02-28"
$ws.Cells.Item(136, 3).Value = "0211"
$ws.Cells.Item(136, 4).Value = "TASHIN"
$ws.Cells.Item(136, 5).Value = 0.23
$ws.Cells.Item(136, 6).Value = 0.23
$ws.Cells.Item(136, 7).Value = 0.21
$ws.Cells.Item(136, 8).Value = 0.21
$ws.Cells.Item(136, 9).Value = 570800

$ws.Range("B135:C136").ClearFormats()
